$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "longitude"
$ws.Range("C1").Value = "latitude"

# City data rows: name, longitude, latitude
$ws.Range("A2").Value = "Tbilisi"
$ws.Range("B2").Value = 44.8271
$ws.Range("C2").Value = 41.7151

$ws.Range("A3").Value = "Batumi"
$ws.Range("B3").Value = 41.6367
$ws.Range("C3").Value = 41.6168

$ws.Range("A4").Value = "Zugdidi"
$ws.Range("B4").Value = 41.8709
$ws.Range("C4").Value = 42.5088

$ws.Range("A5").Value = "Kutaisi"
$ws.Range("B5").Value = 42.718
$ws.Range("C5").Value = 42.2662

$ws.Range("A6").Value = "Gori"
$ws.Range("B6").Value = 44.1083
$ws.Range("C6").Value = 41.9862

# Column widths (character units, matches target file widths as closely as
# the engine's internal pixel-snapped representation allows)
$ws.Range("A:A").ColumnWidth = 7.498697916666667
$ws.Range("B:B").ColumnWidth = 11.166666666666666
$ws.Range("C:C").ColumnWidth = 9.944010416666666

# Match the page setup orientation recorded in the saved file
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("C3").Select()
